# ---------------------------------------------------------------------------
# Add 2022-Q4 data:
#  - insert a new "2022-Q4" worksheet (holding the per-fund breakdown) right
#    after "总计" and before "2022-Q3"
#  - insert a new summary row at the top of the "总计" sheet for 2022-Q4,
#    pushing the existing quarters down by one row
#
# NOTE: worksheet references returned by Worksheets.Item(n) track *tab
# position*, not object identity — once a sheet is inserted/removed the same
# variable can silently start pointing at a different sheet. So we do every
# sheet-order-changing operation FIRST, then re-fetch every sheet we still
# need (by name, which is stable) before touching any cells.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ===========================================================================
# STEP 1 — create the new "2022-Q4" sheet in the right tab position
# ===========================================================================
$q4 = $wb.Worksheets.Add($wb.Worksheets.Item(2))   # inserted before "2022-Q3" -> slot #2
$q4.Name = "2022-Q4"

# Re-fetch everything we need by (now-stable) name.
$total = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Item("2022-Q4")

# ===========================================================================
# STEP 2 — update the "总计" (summary) sheet: insert the 2022-Q4 row at the
# top of the data (row 2), shifting the other quarters down by one row.
# ===========================================================================

# Snapshot the current quarter rows (rows 2..7) before we overwrite anything.
$oldB = @()
$oldC = @()
$oldD = @()
for ($r = 2; $r -le 7; $r++) {
    $oldB += $total.Cells.Item($r, 2).Value()
    $oldC += $total.Cells.Item($r, 3).Value()
    $oldD += $total.Cells.Item($r, 4).Value()
}

# Shift them down into rows 3..8.
for ($i = 0; $i -lt 6; $i++) {
    $r = $i + 3
    $total.Cells.Item($r, 2).Value = $oldB[$i]
    $total.Cells.Item($r, 3).Value = $oldC[$i]
    $total.Cells.Item($r, 4).Value = $oldD[$i]
}

# New row 2 -> the 2022-Q4 summary.
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 10
$total.Cells.Item(2, 4).Value = 1.25

# Row 8's index column (A) needs a value + the same style as the other index
# cells (bold/border, like A2:A7).
$total.Cells.Item(8, 1).Value = 6
$total.Cells.Item(7, 1).Copy()
$total.Cells.Item(8, 1).PasteSpecial(-4122)

# ===========================================================================
# STEP 3 — fill in the new "2022-Q4" sheet with the per-fund breakdown
# ===========================================================================

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")

# column B..H
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q4.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Fund rows: code, name, scale, stock-position, position-ratio, held-value, rank
$rows = @(
    @("016250", "华夏远见成长一年持有混合A", "9.60",  "88.62", "6.36", "0.6106", 2),
    @("016251", "华夏远见成长一年持有混合C", "2.97",  "88.62", "6.36", "0.1889", 2),
    @("005660", "嘉实资源精选股票A",         "2.63",  "93.36", "5.74", "0.1510", 2),
    @("002144", "华安新优选灵活配置混合C",   "15.28", "23.65", "0.73", "0.1115", 10),
    @("001312", "华安新优选灵活配置混合A",   "9.10",  "23.65", "0.73", "0.0664", 10),
    @("005661", "嘉实资源精选股票C",         "1.06",  "93.36", "5.74", "0.0608", 2),
    @("009649", "嘉实精选平衡混合A",         "0.47",  "68.05", "6.32", "0.0297", 3),
    @("011765", "兴银高端制造混合A",         "0.56",  "93.81", "3.11", "0.0174", 4),
    @("011766", "兴银高端制造混合C",         "0.31",  "93.81", "3.11", "0.0096", 4),
    @("009650", "嘉实精选平衡混合C",         "0.04",  "68.05", "6.32", "0.0025", 3)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $fund = $rows[$i]

    $q4.Cells.Item($r, 1).Value = $i               # A: running index (0-based)

    # text-like columns (fund code/name/scale/position/ratio/value) must stay
    # text, even though several of them look like plain numbers.
    $q4.Cells.Item($r, 2).NumberFormat = "@"
    $q4.Cells.Item($r, 2).Value = $fund[0]
    $q4.Cells.Item($r, 3).NumberFormat = "@"
    $q4.Cells.Item($r, 3).Value = $fund[1]
    $q4.Cells.Item($r, 4).NumberFormat = "@"
    $q4.Cells.Item($r, 4).Value = $fund[2]
    $q4.Cells.Item($r, 5).NumberFormat = "@"
    $q4.Cells.Item($r, 5).Value = $fund[3]
    $q4.Cells.Item($r, 6).NumberFormat = "@"
    $q4.Cells.Item($r, 6).Value = $fund[4]
    $q4.Cells.Item($r, 7).NumberFormat = "@"
    $q4.Cells.Item($r, 7).Value = $fund[5]

    $q4.Cells.Item($r, 8).Value = $fund[6]         # H: rank, a real number
}

# Reset the "text trick" styling (NumberFormat "@") back to the plain default
# style used throughout the rest of the workbook, without disturbing the
# values/types we just set.
$cleanCell = $total.Cells.Item(20, 20)             # untouched, default style
$cleanCell.Copy()
$q4.Range("B2:G11").PasteSpecial(-4122)

# Apply the bold/border header style (matches the other quarter sheets).
$q3.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

# Apply the bold/border index-column style (A2:A11), matching A2 on 2022-Q3.
$q3.Cells.Item(2, 1).Copy()
$q4.Range("A2:A11").PasteSpecial(-4122)

$q4.Cells.Item(1, 1).Select()
